# Clean Quarto rebuild for Stats II Spring 2026
#
# Slide 39 ("Main Sources of Spatially Referenced Data" table) had its
# header shortened to "Sources of Spatially Referenced Data". The header
# lives in the first (merged) cell of row 1 of the table that fills the
# slide's single content placeholder.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(39)
$shape = $s.Shapes.Item(1)

$table = $shape.Table
$headerCell = $table.Rows.Item(1).Cells.Item(1)
$headerCell.Shape.TextFrame.TextRange.Text = "Sources of Spatially Referenced Data"
